$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.139.49"
$ws.Range("E2").Value = "  +3.65%  "

$ws.Range("D3").Value = "1.602.09"
$ws.Range("E3").Value = "  +3.55%  "

$ws.Range("E4").Value = "  -0.29%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.84%  "

$ws.Range("E6").Value = "  -0.32%  "

$ws.Range("E7").Value = "  +1.94%  "

$ws.Range("E8").Value = "  +2.47%  "

$ws.Range("E9").Value = "  +1.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.37%  "

$ws.Range("E11").Value = "  +4.81%  "

$ws.Range("D12").Value = "1.825.59"
$ws.Range("E12").Value = "  +3.54%  "

$ws.Range("D13").Value = "1.603.68"
$ws.Range("E13").Value = "  +3.53%  "

$ws.Range("E14").Value = "  +0.43%  "

$ws.Range("E15").Value = "  +1.36%  "

$ws.Range("D16").Value = "26.120.65"
$ws.Range("E16").Value = "  +3.70%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.28%  "

$ws.Range("E18").Value = "  +2.12%  "

$ws.Range("E19").Value = "  -0.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "203.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.88%  "

$ws.Range("E21").Value = "  +3.41%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.70%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.90%  "

$ws.Range("E24").Value = "  +11.57%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.56%  "

$ws.Range("E26").Value = "  -0.27%  "

$ws.Range("E27").Value = "  -4.60%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.56%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.71%  "

$ws.Range("E30").Value = "  +1.78%  "

$ws.Range("E31").Value = "  +2.02%  "

$ws.Range("E32").Value = "  +3.21%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.98"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.79%  "

$ws.Range("E34").Value = "  +1.51%  "

$ws.Range("E35").Value = "  +1.63%  "

$ws.Range("D36").Value = "1.120.76"
$ws.Range("E36").Value = "  +3.45%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0164"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.92%  "

$ws.Range("E39").Value = "  +3.03%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.29"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.15%  "

$ws.Range("E41").Value = "  -0.34%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.779"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.28%  "

$ws.Range("D43").Value = "1.738.24"
$ws.Range("E43").Value = "  +3.54%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.13"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.94%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.27%  "

$ws.Range("E46").Value = "  +3.72%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "53.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.41%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0503"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.32%  "

$ws.Range("E49").Value = "  +1.41%  "

$ws.Range("E50").Value = "  -0.11%  "

$ws.Range("E51").Value = "  -16.69%  "
